$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows (9 and 10) that duplicate rows 7 ("a1") and 8 ("b2"),
# matching the style (s="1") used in column A for the existing label cells.

$ws.Range("A7").Copy($ws.Range("A9"))
$ws.Range("A9").Value = "a1"
$ws.Range("B9").Value = 0.8709480166435242
$ws.Range("C9").Value = 0.8732147216796875
$ws.Range("D9").Value = 269.6632690429688
$ws.Range("E9").Value = 32.62083053588867
$ws.Range("F9").Value = 32.26833343505859
$ws.Range("G9").Value = 181.3432159423828
$ws.Range("H9").Value = 228.9816131591797

$ws.Range("A8").Copy($ws.Range("A10"))
$ws.Range("A10").Value = "b2"
$ws.Range("B10").Value = 0.8709480166435242
$ws.Range("C10").Value = 0.8732147216796875
$ws.Range("D10").Value = 269.6569213867188
$ws.Range("E10").Value = 32.6202278137207
$ws.Range("F10").Value = 32.26775360107422
$ws.Range("G10").Value = 181.3432006835938
$ws.Range("H10").Value = 228.9815063476562
